$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns with refreshed
# crypto data. D-column values that look like plain numbers are first
# forced to Text format so values such as "0.998" or "12.09" are kept
# exactly as text (matching the source data) instead of being
# auto-converted to numbers and losing significant trailing zeros.

$ws.Range("D2").Value = '55.531.24'
$ws.Range("E2").Value = '  -3.10%  '

$ws.Range("D3").Value = '2.967.86'
$ws.Range("E3").Value = '  -5.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '488.75'
$ws.Range("E5").Value = '  -5.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.77'
$ws.Range("E6").Value = '  -0.92%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("D8").Value = '2.972.79'
$ws.Range("E8").Value = '  -5.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.417'
$ws.Range("E9").Value = '  -5.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.09'
$ws.Range("E10").Value = '  -1.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.101'
$ws.Range("E11").Value = '  -6.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.349'
$ws.Range("E12").Value = '  -8.12%  '

$ws.Range("E13").Value = '  +0.34%  '

$ws.Range("D14").Value = '3.472.15'
$ws.Range("E14").Value = '  -5.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '24.76'
$ws.Range("E15").Value = '  -1.67%  '

$ws.Range("D16").Value = '55.400.29'
$ws.Range("E16").Value = '  -3.41%  '

$ws.Range("D17").Value = '2.952.38'
$ws.Range("E17").Value = '  -5.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000140'
$ws.Range("E18").Value = '  -5.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.58'
$ws.Range("E19").Value = '  -2.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.09'
$ws.Range("E20").Value = '  -5.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.46'
$ws.Range("E21").Value = '  -5.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '319.44'
$ws.Range("E22").Value = '  -6.41%  '

$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.463'
$ws.Range("E24").Value = '  -7.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '60.27'
$ws.Range("E25").Value = '  -11.82%  '

$ws.Range("E26").Value = '  -0.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.163'
$ws.Range("E27").Value = '  +0.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'

$ws.Range("D29").Value = '0.0₃0844'
$ws.Range("E29").Value = '  -8.65%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.49'
$ws.Range("E30").Value = '  -2.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.52'
$ws.Range("E31").Value = '  -3.92%  '

$ws.Range("E32").Value = '  -1.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.69'
$ws.Range("E33").Value = '  -7.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.40'
$ws.Range("E34").Value = '  -9.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '149.10'
$ws.Range("E35").Value = '  -5.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.38'
$ws.Range("E36").Value = '  -8.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.29'
$ws.Range("E37").Value = '  -5.99%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.67'
$ws.Range("E38").Value = '  -7.66%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.25'
$ws.Range("E39").Value = '  -9.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0647'
$ws.Range("E40").Value = '  -4.43%  '

$ws.Range("D41").Value = '2.993.64'
$ws.Range("E41").Value = '  -5.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("E42").Value = '  -0.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '36.03'
$ws.Range("E43").Value = '  -10.42%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  -5.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.630'
$ws.Range("E45").Value = '  -8.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.37'
$ws.Range("E46").Value = '  -5.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.54'
$ws.Range("E47").Value = '  -7.85%  '

$ws.Range("D48").Value = '2.128.27'
$ws.Range("E48").Value = '  -4.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0233'
$ws.Range("E49").Value = '  +0.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.11'
$ws.Range("E50").Value = '  -3.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.57'
$ws.Range("E51").Value = '  -8.38%  '
